# The author switched the presentation's applied Design/Theme away from the
# custom "Integral" (Red Violet) theme back to the default, built-in
# "Office Theme" palette (Design > Themes > Office Theme in the PowerPoint
# ribbon). That swaps the 12 theme colours used across every slide, slide
# master and layout (the font scheme and format/effect scheme of the
# built-in "Office Theme" are the same ones already used by this deck, so
# only the colour scheme actually changes).

$p = $ppt.ActivePresentation

# Helper: build a COM/OLE RGB long (0x00BBGGRR) from a "RRGGBB" hex string.
function HexToRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme colour scheme, in ThemeColorScheme slot order:
# 1 Dark1, 2 Light1, 3 Dark2, 4 Light2,
# 5-10 Accent1..Accent6, 11 Hyperlink, 12 FollowedHyperlink
$officeThemeColors = @(
    "000000", # Dark 1
    "FFFFFF", # Light 1
    "44546A", # Dark 2
    "E7E6E6", # Light 2
    "5B9BD5", # Accent 1
    "ED7D31", # Accent 2
    "A5A5A5", # Accent 3
    "FFC000", # Accent 4
    "4472C4", # Accent 5
    "70AD47", # Accent 6
    "0563C1", # Hyperlink
    "954F72"  # Followed Hyperlink
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $themeColors.Item($i).RGB = HexToRgb($officeThemeColors[$i - 1])
}
